$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole data block to Text format so numeric-looking values
# (prices, hour) are written back as text, matching the inlineStr cells
# already used throughout this sheet, then restore the default "Normal"
# style so no new per-cell formatting is introduced.
$dataRange = $ws.Range("B2:G51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '246.16'
$ws.Range("G2").Value = '23'
$ws.Range("D3").Value = '24.05'
$ws.Range("G3").Value = '23'
$ws.Range("D4").Value = '5.357'
$ws.Range("G4").Value = '23'
$ws.Range("D5").Value = '0.05811'
$ws.Range("G5").Value = '23'
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '6.477'
$ws.Range("E6").Value = '5KuCoinTokenKCS'
$ws.Range("G6").Value = '23'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '3.366'
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("G7").Value = '23'
$ws.Range("D8").Value = '0.8101'
$ws.Range("G8").Value = '23'
$ws.Range("D9").Value = '0.9198'
$ws.Range("G9").Value = '23'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1404'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").Value = '23'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07397'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("G11").Value = '23'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03208'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G12").Value = '23'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03032'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("G13").Value = '23'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09379'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("G14").Value = '23'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '3.848'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").Value = '23'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001567'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("G16").Value = '23'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04698'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("G17").Value = '23'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '0.0005987'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("G18").Value = '23'
$ws.Range("D19").Value = '0.006080'
$ws.Range("G19").Value = '23'
$ws.Range("D20").Value = '0.001247'
$ws.Range("G20").Value = '23'
$ws.Range("G21").Value = '23'
$ws.Range("D22").Value = '0.00008796'
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("G22").Value = '23'
$ws.Range("D23").Value = '3.597'
$ws.Range("G23").Value = '23'
$ws.Range("G24").Value = '23'
$ws.Range("D25").Value = '0.3179'
$ws.Range("G25").Value = '23'
$ws.Range("D26").Value = '0.1318'
$ws.Range("G26").Value = '23'
$ws.Range("G27").Value = '23'
$ws.Range("G28").Value = '23'
$ws.Range("G29").Value = '23'
$ws.Range("G30").Value = '23'
$ws.Range("G31").Value = '23'
$ws.Range("G32").Value = '23'
$ws.Range("G33").Value = '23'
$ws.Range("G34").Value = '23'
$ws.Range("G35").Value = '23'
$ws.Range("G36").Value = '23'
$ws.Range("G37").Value = '23'
$ws.Range("G38").Value = '23'
$ws.Range("G39").Value = '23'
$ws.Range("D40").Value = '0.03844'
$ws.Range("G40").Value = '23'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1068'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("G41").Value = '23'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '0.002749'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("G42").Value = '23'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '0.003040'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("G43").Value = '23'
$ws.Range("D44").Value = '0.009066'
$ws.Range("G44").Value = '23'
$ws.Range("D45").Value = '0.00005247'
$ws.Range("G45").Value = '23'
$ws.Range("G46").Value = '23'
$ws.Range("D47").Value = '0.7096'
$ws.Range("G47").Value = '23'
$ws.Range("D48").Value = '0.001832'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("G48").Value = '23'
$ws.Range("G49").Value = '23'
$ws.Range("D50").Value = '0.0001999'
$ws.Range("G50").Value = '23'
$ws.Range("G51").Value = '23'

# Restore the default style on the whole block (NumberFormat alone leaves
# a transient "@" text-format style reference on touched cells).
$dataRange.Style = "Normal"
